$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 29425
$ws.Range("J63").Value = 29425
$ws.Range("L63").Value = 29425
$ws.Range("N63").Value = -30673

$ws.Range("H66").Value = 29425
$ws.Range("J66").Value = 29425
$ws.Range("L66").Value = 88275
$ws.Range("N66").Value = -94515

$ws.Range("H76").Value = 2692.3635
$ws.Range("I76").Value = 2648.9048
$ws.Range("K76").Value = 2648.9048
$ws.Range("M76").Value = -2333.9048

$ws.Range("H79").Value = 2692.3635
$ws.Range("I79").Value = 2648.9048
$ws.Range("K79").Value = 2648.9048
$ws.Range("M79").Value = -1556.9048

$ws.Range("H88").Value = 9506.933999999999
$ws.Range("I88").Value = 7232.5
$ws.Range("J88").Value = 12106.286
$ws.Range("K88").Value = 7232.5
$ws.Range("L88").Value = 12106.286
$ws.Range("M88").Value = -6826.5
$ws.Range("N88").Value = -12918.286

$ws.Range("H91").Value = 9506.933999999999
$ws.Range("I91").Value = 7232.5
$ws.Range("J91").Value = 12106.286
$ws.Range("K91").Value = 7232.5
$ws.Range("L91").Value = 12106.286
$ws.Range("M91").Value = -5828.5
$ws.Range("N91").Value = -14914.286

$ws.Range("H106").Value = 2727.1333
$ws.Range("I106").Value = 1901.2858
$ws.Range("K106").Value = 1901.2858
$ws.Range("M106").Value = -1270.2858

$ws.Range("H107").Value = 664.6429000000001
$ws.Range("I107").Value = 600.55554
$ws.Range("K107").Value = 600.55554
$ws.Range("M107").Value = 1319.44446

$ws.Range("H132").Value = 3848371.5
$ws.Range("I132").Value = 4446485.5
$ws.Range("K132").Value = 13339456.5
$ws.Range("M132").Value = -13336926.5

$ws.Range("H137").Value = 2642.106
$ws.Range("J137").Value = 2688.647
$ws.Range("L137").Value = 8065.941
$ws.Range("N137").Value = -13165.941

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6634.1836
$ws.Range("I32").Value = 5301.722
$ws.Range("J32").Value = 21624.375
$ws.Range("K32").Value = 5301.722
$ws.Range("L32").Value = 21624.375
$ws.Range("M32").Value = -5014.722
$ws.Range("N32").Value = -22198.375

$ws.Range("H45").Value = 1492.7833
$ws.Range("I45").Value = 1127.6522
$ws.Range("J45").Value = 2692.5
$ws.Range("K45").Value = 1127.6522
$ws.Range("L45").Value = 2692.5
$ws.Range("M45").Value = -750.6522
$ws.Range("N45").Value = -3446.5

$ws.Range("H63").Value = 1953.1333
$ws.Range("I63").Value = 1935.5
$ws.Range("K63").Value = 1935.5
$ws.Range("M63").Value = -1249.5

$ws.Range("H66").Value = 1953.1333
$ws.Range("I66").Value = 1935.5
$ws.Range("K66").Value = 9677.5
$ws.Range("M66").Value = -6245.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 30647.5
$ws.Range("J69").Value = 30647.5
$ws.Range("L69").Value = 30647.5
$ws.Range("N69").Value = -32269.5

$ws.Range("H72").Value = 30647.5
$ws.Range("J72").Value = 30647.5
$ws.Range("L72").Value = 91942.5
$ws.Range("N72").Value = -100054.5

$ws.Range("H105").Value = 1633.5
$ws.Range("I105").Value = 1519.25
$ws.Range("K105").Value = 1519.25
$ws.Range("M105").Value = 227.75

$ws.Range("H125").Value = 24000
$ws.Range("J125").Value = 24000
$ws.Range("L125").Value = 24000
$ws.Range("N125").Value = -33840

$ws.Range("H134").Value = 3246.389
$ws.Range("I134").Value = 2032.3846
$ws.Range("J134").Value = 6402.8
$ws.Range("K134").Value = 6097.1538
$ws.Range("L134").Value = 19208.4
$ws.Range("M134").Value = -3562.1538
$ws.Range("N134").Value = -24278.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3194
$ws.Range("I99").Value = 1490
$ws.Range("J99").Value = 5750
$ws.Range("K99").Value = 1490
$ws.Range("L99").Value = 5750
$ws.Range("M99").Value = 8
$ws.Range("N99").Value = -8746

$ws.Range("H126").Value = 3194
$ws.Range("I126").Value = 1490
$ws.Range("J126").Value = 5750
$ws.Range("K126").Value = 4470
$ws.Range("L126").Value = 17250
$ws.Range("M126").Value = -2000
$ws.Range("N126").Value = -22190

$ws.Range("H134").Value = 2476.5
$ws.Range("I134").Value = 829.2727
$ws.Range("K134").Value = 2487.8181
$ws.Range("M134").Value = 47.18190000000004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3275.5557
$ws.Range("I134").Value = 2702
$ws.Range("J134").Value = 3992.5
$ws.Range("K134").Value = 8106
$ws.Range("L134").Value = 11977.5
$ws.Range("M134").Value = -3036
$ws.Range("N134").Value = -22117.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5014.2144
$ws.Range("I70").Value = 4742.7144
$ws.Range("J70").Value = 5285.7144
$ws.Range("K70").Value = 4742.7144
$ws.Range("L70").Value = 5285.7144
$ws.Range("M70").Value = -4472.7144
$ws.Range("N70").Value = -5825.7144

$ws.Range("H73").Value = 5014.2144
$ws.Range("I73").Value = 4742.7144
$ws.Range("J73").Value = 5285.7144
$ws.Range("K73").Value = 4742.7144
$ws.Range("L73").Value = 5285.7144
$ws.Range("M73").Value = -3806.7144
$ws.Range("N73").Value = -7157.7144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 902.5714
$ws.Range("I16").Value = 1011
$ws.Range("J16").Value = 252
$ws.Range("K16").Value = 1011
$ws.Range("L16").Value = 252
$ws.Range("M16").Value = -841
$ws.Range("N16").Value = -592

$ws.Range("H125").Value = 29833.334
$ws.Range("J125").Value = 29833.334
$ws.Range("L125").Value = 29833.334
$ws.Range("N125").Value = -39673.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3450465
$ws.Range("I126").Value = 1651.4783
$ws.Range("J126").Value = 16670917
$ws.Range("K126").Value = 4954.4349
$ws.Range("L126").Value = 50012751
$ws.Range("M126").Value = -2484.4349
$ws.Range("N126").Value = -50017691

$ws.Range("H132").Value = 3794.0908
$ws.Range("I132").Value = 1578.2667
$ws.Range("J132").Value = 13765.3
$ws.Range("K132").Value = 4734.800099999999
$ws.Range("L132").Value = 41295.89999999999
$ws.Range("M132").Value = -2204.800099999999
$ws.Range("N132").Value = -46355.89999999999

$ws.Range("H136").Value = 1074.2916
$ws.Range("J136").Value = 1805.5555
$ws.Range("L136").Value = 5416.666499999999
$ws.Range("N136").Value = -10516.6665
